$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.7533
$ws.Range("E2").Value = 0.7151999999999999
$ws.Range("F2").Value = 0.7958
$ws.Range("G2").Value = 0.7533
